$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: "BioSum 5.9.0 Release Notes" / "3 May, 2022"
#    -> "BioSum 5.10.0 Release Notes" / "January 13, 2023"
#    Also relocate the lone "_GoBack" bookmark from the end of the document to
#    the very start of the title paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
# Insert a throwaway placeholder character at position 0 so that a bookmark
# can be collapsed immediately after it (collapsing directly at position 0
# expands to the next paragraph in this runtime), then remove the
# placeholder once the bookmark is anchored.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()

$d.Content.Find.Execute("BioSum 5.9.0 Release Notes", $false, $false, $false, $false, $false, $true, 1, $false, "BioSum 5.10.0 Release Notes", 2)
$d.Content.Find.Execute("3 May, 2022", $false, $false, $false, $false, $false, $true, 1, $false, "January 13, 2023", 2)

# ---------------------------------------------------------------------------
# 2. "Current Release Notes (5.9.0)" -> "Current Release Notes (5.10.0)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Current Release Notes (5.9.0)", $false, $false, $false, $false, $false, $true, 1, $false, "Current Release Notes (5.10.0)", 2)

# ---------------------------------------------------------------------------
# 3. First bullet: "New FICS (...) ..." -> "FVS_CutTree table has been
#    converted to SQLite: ..."
# ---------------------------------------------------------------------------
$ficsIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("New FICS")) {
        $ficsIdx = $i
        break
    }
}
$pFics = $d.Paragraphs.Item($ficsIdx)
$r = $pFics.Range
$r.MoveEnd(1, -1)
$r.Text = ""

$boldText = "FVS_CutTree table has been converted to SQLite"
$restText = ": the FVS_CutTree (aka BiosumCalc FVS_Tree) table is now written to the /fvs/data/FVSOUT_TREE_LIST.db. This means that legacy projects are frozen at v5.9.0. Optimizer scenarios can be added, changed, and executed, but Processor can no longer run on these older projects because they don't have an FVSOUT_TREE_LIST.db. It may be possible to manually load the FVS_CutTree table from each FVS_Tree table in the BiosumCalc director(ies) if necessary. Please contact BioSum support for instructions."

$pFics.Range.InsertAfter($boldText + $restText)

$full = $pFics.Range
$full.MoveEnd(1, -1)
$boldRange = $full.Duplicate
$boldRange.SetRange($full.Start, $full.Start + $boldText.Length)
$boldRange.Font.Bold = 1

$endPt = $pFics.Range
$endPt.MoveEnd(1, -1)
$endPt.Collapse(0)
$endPt.InsertAfter([char]11)

# ---------------------------------------------------------------------------
# 4. Delete the two paragraphs that used to follow the FICS bullet:
#    "Important note for previous Forest Service users: ..." and
#    "Note for all previous BioSum users: ..."
# ---------------------------------------------------------------------------
$pImportant = $d.Paragraphs.Item($ficsIdx + 1)
$pNoteAll = $d.Paragraphs.Item($ficsIdx + 2)
$delRange = $d.Range($pImportant.Range.Start, $pNoteAll.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 5. Bullet "Integration with FVSOn*: ..." -> "Updates to FVS sequence
#    number definition screen: ..." and append a trailing line break.
# ---------------------------------------------------------------------------
$fvsOnIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Integration with")) {
        $fvsOnIdx = $i
        break
    }
}
$pFvsOn = $d.Paragraphs.Item($fvsOnIdx)
$r2 = $pFvsOn.Range
$r2.MoveEnd(1, -1)
$r2.Text = ""

$boldText2 = "Updates to FVS sequence number definition screen"
$restText2 = ": Functional and cosmetic changes have been made to the FVS sequence number definition screen to make it easier to use and understand. Sequence numbers can be assigned to multiple tables at the same time. The selection of sequence number templates has been updated to reflect commonly used BioSum configurations."

$pFvsOn.Range.InsertAfter($boldText2 + $restText2)

$full2 = $pFvsOn.Range
$full2.MoveEnd(1, -1)
$boldRange2 = $full2.Duplicate
$boldRange2.SetRange($full2.Start, $full2.Start + $boldText2.Length)
$boldRange2.Font.Bold = 1

$endPt2 = $pFvsOn.Range
$endPt2.MoveEnd(1, -1)
$endPt2.Collapse(0)
$endPt2.InsertAfter([char]11)

# ---------------------------------------------------------------------------
# 6. Delete the paragraph that used to follow it:
#    "FVSOn* generates output only in SQLite, ..."
# ---------------------------------------------------------------------------
$pFvsOnGen = $d.Paragraphs.Item($fvsOnIdx + 1)
$pFvsOnGen.Range.Delete()

# ---------------------------------------------------------------------------
# 7. Delete bullet "Cumulative SQLite POP tables: ..." entirely.
# ---------------------------------------------------------------------------
$popIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Cumulative SQLite POP tables")) {
        $popIdx = $i
        break
    }
}
$pPop = $d.Paragraphs.Item($popIdx)
$pPop.Range.Delete()

# ---------------------------------------------------------------------------
# 8. Final bullet: "Several additional minor enhancements ..." ->
#    "Many additional minor enhancements are described in the BioSum
#    online release documentation." with a real hyperlink.
# ---------------------------------------------------------------------------
$lastIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Several additional")) {
        $lastIdx = $i
        break
    }
}
$pLast = $d.Paragraphs.Item($lastIdx)
$r3 = $pLast.Range
$r3.MoveEnd(1, -1)
$r3.Text = ""

$pLast.Range.InsertAfter("Many additional minor enhancements are described in the BioSum ")
$afterLead = $pLast.Range
$afterLead.MoveEnd(1, -1)
$afterLead.Collapse(0)

$linkStart = $afterLead.Start
$linkDisplay = "online release documentation"
$afterLead.InsertAfter($linkDisplay)
$linkEnd = $linkStart + $linkDisplay.Length
$linkRange = $d.Range($linkStart, $linkEnd)
$d.Hyperlinks.Add($linkRange, "https://github.com/USFS-PNW/Fia-Biosum-Manager/releases/tag/5.10.0", $null, $null, $linkDisplay)

$pLastAgain = $d.Paragraphs.Item($lastIdx)
$tailPt = $pLastAgain.Range
$tailPt.MoveEnd(1, -1)
$tailPt.Collapse(0)
$tailPt.InsertAfter(".")
